# Insert a new "time" column before the existing "venue_id" column.
# Template header row goes from:
#   A1=event_id | B1=event_name | C1=date | D1=venue_id
# to:
#   A1=event_id | B1=event_name | C1=date | D1=time | E1=venue_id

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the old D1 ("venue_id") content one column to the right, into E1.
$ws.Range("E1").Value = $ws.Range("D1").Value()

# Put the new header in the now-vacated D1.
$ws.Range("D1").Value = "time"

# Leave the selection where the author left it.
$ws.Range("D4").Select()
